$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the original values for the rows that are being cyclically
# rotated (row2 <- row3, row3 <- row4, row4 <- row2) for the columns
# that vary: D (Fecha), M (Volumen), N (Precio minimo), O (Precio maximo),
# P (Precio promedio ponderado), R (Origen), S (Precio $/Kg).

$cols = @("D", "M", "N", "O", "P", "R", "S")

$orig2 = @{}
$orig3 = @{}
$orig4 = @{}
foreach ($col in $cols) {
    $orig2[$col] = $ws.Range("${col}2").Value2
    $orig3[$col] = $ws.Range("${col}3").Value2
    $orig4[$col] = $ws.Range("${col}4").Value2
}

foreach ($col in $cols) {
    $ws.Range("${col}2").Value2 = $orig3[$col]
    $ws.Range("${col}3").Value2 = $orig4[$col]
    $ws.Range("${col}4").Value2 = $orig2[$col]
}
